# Updates the cryptos price/volume table with the latest scraped values.
# Note: Price (column D) values are written with a leading apostrophe so
# Excel stores them as literal text (quote-prefixed) instead of silently
# coercing decimal-looking strings (e.g. "416.60") into floating point
# numbers, which would introduce binary rounding drift when saved back
# to OOXML (e.g. "416.60" -> 416.60000000000002).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''63.116.41'
$ws.Range('E2').Value = '  +1.91%  '
$ws.Range('D3').Value = '''3.488.36'
$ws.Range('E3').Value = '  +2.33%  '
$ws.Range('E4').Value = '  -0.22%  '
$ws.Range('D5').Value = '''416.60'
$ws.Range('E5').Value = '  +1.71%  '
$ws.Range('D6').Value = '''130.43'
$ws.Range('E6').Value = '  +1.29%  '
$ws.Range('E7').Value = '  -1.33%  '
$ws.Range('D8').Value = '''0.999'
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').Value = '''0.736'
$ws.Range('E9').Value = '  +0.27%  '
$ws.Range('D10').Value = '''0.155'
$ws.Range('E10').Value = '  +9.83%  '
$ws.Range('D11').Value = '''42.64'
$ws.Range('E11').Value = '  -1.16%  '
$ws.Range('D12').Value = '''9.82'
$ws.Range('E12').Value = '  +5.12%  '
$ws.Range('E13').Value = '  +4.73%  '
$ws.Range('D14').Value = '''4.035.74'
$ws.Range('E14').Value = '  +2.15%  '
$ws.Range('E15').Value = '  -0.18%  '
$ws.Range('D16').Value = '''20.63'
$ws.Range('E16').Value = '  -3.59%  '
$ws.Range('D17').Value = '''3.494.91'
$ws.Range('E17').Value = '  +2.60%  '
$ws.Range('D18').Value = '''12.70'
$ws.Range('E18').Value = '  +1.68%  '
$ws.Range('D19').Value = '''1.09'
$ws.Range('E19').Value = '  -0.32%  '
$ws.Range('D20').Value = '''62.963.07'
$ws.Range('E20').Value = '  +1.64%  '
$ws.Range('D21').Value = '''467.43'
$ws.Range('E21').Value = '  +3.96%  '
$ws.Range('D22').Value = '''90.84'
$ws.Range('E22').Value = '  -0.64%  '
$ws.Range('D23').Value = '''3.31'
$ws.Range('E23').Value = '  +3.28%  '
$ws.Range('E24').Value = '  +0.73%  '
$ws.Range('D25').Value = '''10.73'
$ws.Range('E25').Value = '  +14.14%  '
$ws.Range('D26').Value = '''3.33'
$ws.Range('E26').Value = '  +0.89%  '
$ws.Range('D27').Value = '''33.76'
$ws.Range('E27').Value = '  +2.03%  '
$ws.Range('D28').Value = '''4.86'
$ws.Range('E28').Value = '  +1.30%  '
$ws.Range('D29').Value = '''7.57'
$ws.Range('E29').Value = '  -1.35%  '
$ws.Range('E30').Value = '  +1.17%  '
$ws.Range('E31').Value = '  -0.75%  '
$ws.Range('D32').Value = '''0.169'
$ws.Range('E33').Value = '  -0.90%  '
$ws.Range('D34').Value = '''41.12'
$ws.Range('E34').Value = '  -3.43%  '
$ws.Range('E35').Value = '  +0.18%  '
$ws.Range('D36').Value = '''58.25'
$ws.Range('E36').Value = '  +8.20%  '
$ws.Range('D37').Value = '''0.0492'
$ws.Range('E37').Value = '  -2.66%  '
$ws.Range('D38').Value = '''0.998'
$ws.Range('E38').Value = '  +0.00%  '
$ws.Range('D39').Value = '''3.06'
$ws.Range('E39').Value = '  +3.44%  '
$ws.Range('D40').Value = '''2.75'
$ws.Range('E40').Value = '  +7.92%  '
$ws.Range('B41').Value = 'Monero'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D41').Value = '''148.32'
$ws.Range('E41').Value = '  +2.97%  '
$ws.Range('E42').Value = '  -0.78%  '
$ws.Range('B43').Value = 'TheGraph'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D43').Value = '''0.323'
$ws.Range('E43').Value = '  +1.21%  '
$ws.Range('D44').Value = '''3.34'
$ws.Range('E44').Value = '  -1.77%  '
$ws.Range('D45').Value = '''4.44'
$ws.Range('E45').Value = '  +2.56%  '
$ws.Range('D46').Value = '''2.07'
$ws.Range('E46').Value = '  +2.78%  '
$ws.Range('D47').Value = '''0.0₃0591'
$ws.Range('E47').Value = '  +39.46%  '
$ws.Range('E48').Value = '  +10.14%  '
$ws.Range('D49').Value = '''16.43'
$ws.Range('E49').Value = '  -1.29%  '
$ws.Range('D50').Value = '''22.23'
$ws.Range('E50').Value = '  -1.52%  '
$ws.Range('E51').Value = '  -5.36%  '
